$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -15
$ws.Range("F4").Value = -12
$ws.Range("F5").Value = 13
$ws.Range("F8").Value = -10
